$wb = $excel.ActiveWorkbook

# Rename sheets from *img to img* (e.g. "himg" -> "imgh")
$wb.Worksheets.Item("himg").Name = "imgh"
$wb.Worksheets.Item("timg").Name = "imgt"
$wb.Worksheets.Item("simg").Name = "imgs"
$wb.Worksheets.Item("gimg").Name = "imgg"
$wb.Worksheets.Item("wimg").Name = "imgw"
$wb.Worksheets.Item("bimg").Name = "imgb"
$wb.Worksheets.Item("eimg").Name = "imge"

# Make the last sheet (imge) the active/selected tab
$wb.Worksheets.Item("imge").Activate()
